$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each crypto row's Price (D) and Volume(1h) (E) columns are stored as literal
# text (e.g. "29.707.24", "  -2.93%  "), never as real numbers/percentages.
# Force NumberFormat to Text before assigning so Excel does not reinterpret
# numeric-looking strings (like "345.00" or "1.009") as actual numbers, then
# restore the Normal style so no stray formatting is left behind.
function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.707.24"
Set-TextValue "E2" "  -2.93%  "
Set-TextValue "D3" "2.095.89"
Set-TextValue "E3" "  -2.12%  "
Set-TextValue "D4" "1.009"
Set-TextValue "E4" "  +0.15%  "
Set-TextValue "D5" "345.00"
Set-TextValue "E5" "  -2.14%  "
Set-TextValue "E6" "  +0.10%  "
Set-TextValue "D7" "0.5150"
Set-TextValue "E7" "  -2.16%  "
Set-TextValue "D8" "0.4393"
Set-TextValue "E8" "  -3.58%  "
Set-TextValue "D9" "52.50"
Set-TextValue "E9" "  -2.23%  "
Set-TextValue "D10" "0.09253"
Set-TextValue "E10" "  +1.23%  "
Set-TextValue "E11" "  -1.08%  "
Set-TextValue "D12" "24.87"
Set-TextValue "E12" "  -0.27%  "
Set-TextValue "D13" "2.088.15"
Set-TextValue "E13" "  -2.52%  "
Set-TextValue "D14" "8.283"
Set-TextValue "E14" "  +1.65%  "
Set-TextValue "D15" "6.743"
Set-TextValue "E15" "  -1.94%  "
Set-TextValue "D16" "99.41"
Set-TextValue "E16" "  -2.71%  "
Set-TextValue "E17" "  -1.67%  "
Set-TextValue "E18" "  +0.12%  "
Set-TextValue "D19" "20.86"
Set-TextValue "E19" "  +6.12%  "
Set-TextValue "D20" "0.06664"
Set-TextValue "E20" "  -0.58%  "
Set-TextValue "E21" "  +0.05%  "
Set-TextValue "D22" "6.193"
Set-TextValue "E22" "  -2.44%  "
Set-TextValue "D23" "29.746.35"
Set-TextValue "E23" "  -3.13%  "
Set-TextValue "E24" "  -2.18%  "
Set-TextValue "D25" "2.320"
Set-TextValue "E25" "  -2.87%  "
Set-TextValue "D26" "2.343.13"
Set-TextValue "E26" "  -1.07%  "
Set-TextValue "E27" "  -2.62%  "
Set-TextValue "D28" "2.524"
Set-TextValue "E28" "  -4.77%  "
Set-TextValue "D29" "161.94"
Set-TextValue "E29" "  -1.64%  "
Set-TextValue "D30" "133.12"
Set-TextValue "E30" "  -2.42%  "
Set-TextValue "D31" "1.133"
Set-TextValue "E31" "  -7.30%  "
Set-TextValue "E32" "  -2.79%  "
Set-TextValue "D33" "1.651"
Set-TextValue "E33" "  -1.18%  "
Set-TextValue "D34" "6.173"
Set-TextValue "E34" "  -3.20%  "
Set-TextValue "D35" "3.936"
Set-TextValue "E35" "  -1.85%  "
Set-TextValue "D36" "6.189"
Set-TextValue "E36" "  +0.21%  "
Set-TextValue "D37" "10.28"
Set-TextValue "E37" "  -1.94%  "
Set-TextValue "D39" "0.06703"
Set-TextValue "E39" "  -3.30%  "
Set-TextValue "E40" "  -1.69%  "
Set-TextValue "D41" "0.6857"
Set-TextValue "E41" "  -2.16%  "
Set-TextValue "D42" "0.2225"
Set-TextValue "E42" "  -4.83%  "
Set-TextValue "D43" "1.297"
Set-TextValue "E43" "  +2.05%  "
Set-TextValue "D44" "0.6628"
Set-TextValue "E44" "  +2.55%  "
Set-TextValue "D45" "14.27"
Set-TextValue "E45" "  -3.78%  "
Set-TextValue "D46" "2.315"
Set-TextValue "E46" "  -1.70%  "
Set-TextValue "E47" "  -3.39%  "
Set-TextValue "E48" "  -5.98%  "
Set-TextValue "D49" "1.220"
Set-TextValue "D50" "82.23"
Set-TextValue "E50" "  -0.94%  "
Set-TextValue "D51" "0.3306"
Set-TextValue "E51" "  +0.75%  "
